# This workbook is a NATMI ligand-receptor edge table for Sema3a -> Nrp1
# across the ECs / FAPs / MuSCs clusters. The commit "update scripts wuth
# new tpm" re-ran the analysis pipeline against a refreshed TPM expression
# matrix. That changed the raw Sema3a ligand statistics for the "ECs"
# sending cluster (columns E-H) and the raw Nrp1 receptor statistics for the
# "ECs" target cluster (columns M-N), which in turn changes every derived
# specificity / edge-weight column for all nine sending x target rows:
#   I/J = ligand value (G/H) as a share of the total across sending clusters
#   O/P = receptor value (M/N) as a share of the total across target clusters
#   Q/R = ligand value (G/H) x receptor value (M/N)   -> edge weight
#   S/T = ligand share (I/J) x receptor share (O/P)   -> edge specificity
#
# Values below are the regenerated figures from the updated pipeline run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.250631
$ws.Range("H2").Value = 0.751893
$ws.Range("I2").Value = 0.2648339568266264
$ws.Range("J2").Value = 0.2648339568266264
$ws.Range("M2").Value = 123.2806423333333
$ws.Range("N2").Value = 369.841927
$ws.Range("O2").Value = 0.6241574062367528
$ws.Range("P2").Value = 0.6241574062367526
$ws.Range("Q2").Value = 30.89795066864567
$ws.Range("R2").Value = 278.081556017811
$ws.Range("S2").Value = 0.1652980755763233
$ws.Range("T2").Value = 0.1652980755763232
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.250631
$ws.Range("H3").Value = 0.751893
$ws.Range("I3").Value = 0.2648339568266264
$ws.Range("J3").Value = 0.2648339568266264
$ws.Range("O3").Value = 0.2392728888301323
$ws.Range("P3").Value = 0.2392728888301322
$ws.Range("Q3").Value = 11.84483568014166
$ws.Range("R3").Value = 106.603521121275
$ws.Range("S3").Value = 0.06336758591022142
$ws.Range("T3").Value = 0.0633675859102214
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.250631
$ws.Range("H4").Value = 0.751893
$ws.Range("I4").Value = 0.2648339568266264
$ws.Range("J4").Value = 0.2648339568266264
$ws.Range("O4").Value = 0.136569704933115
$ws.Range("P4").Value = 0.136569704933115
$ws.Range("Q4").Value = 6.760672810560667
$ws.Range("R4").Value = 60.84605529504601
$ws.Range("S4").Value = 0.03616829534008169
$ws.Range("T4").Value = 0.03616829534008168
# Row 5
$ws.Range("I5").Value = 0.2480790641859371
$ws.Range("J5").Value = 0.2480790641859371
$ws.Range("M5").Value = 123.2806423333333
$ws.Range("N5").Value = 369.841927
$ws.Range("O5").Value = 0.6241574062367528
$ws.Range("P5").Value = 0.6241574062367526
$ws.Range("Q5").Value = 28.94317171026089
$ws.Range("R5").Value = 260.488545392348
$ws.Range("S5").Value = 0.1548403852439354
$ws.Range("T5").Value = 0.1548403852439354
# Row 6
$ws.Range("I6").Value = 0.2480790641859371
$ws.Range("J6").Value = 0.2480790641859371
$ws.Range("O6").Value = 0.2392728888301323
$ws.Range("P6").Value = 0.2392728888301322
$ws.Range("S6").Value = 0.05935859434604498
$ws.Range("T6").Value = 0.05935859434604497
# Row 7
$ws.Range("I7").Value = 0.2480790641859371
$ws.Range("J7").Value = 0.2480790641859371
$ws.Range("O7").Value = 0.136569704933115
$ws.Range("P7").Value = 0.136569704933115
$ws.Range("S7").Value = 0.03388008459595674
$ws.Range("T7").Value = 0.03388008459595673
# Row 8
$ws.Range("G8").Value = 0.4609646666666666
$ws.Range("I8").Value = 0.4870869789874365
$ws.Range("J8").Value = 0.4870869789874365
$ws.Range("M8").Value = 123.2806423333333
$ws.Range("N8").Value = 369.841927
$ws.Range("O8").Value = 0.6241574062367528
$ws.Range("P8").Value = 0.6241574062367526
$ws.Range("Q8").Value = 56.82802019963755
$ws.Range("R8").Value = 511.4521817967379
$ws.Range("S8").Value = 0.304018945416494
$ws.Range("T8").Value = 0.304018945416494
# Row 9
$ws.Range("G9").Value = 0.4609646666666666
$ws.Range("I9").Value = 0.4870869789874365
$ws.Range("J9").Value = 0.4870869789874365
$ws.Range("O9").Value = 0.2392728888301323
$ws.Range("P9").Value = 0.2392728888301322
$ws.Range("R9").Value = 196.0669533264499
$ws.Range("S9").Value = 0.1165467085738659
$ws.Range("T9").Value = 0.1165467085738659
# Row 10
$ws.Range("G10").Value = 0.4609646666666666
$ws.Range("I10").Value = 0.4870869789874365
$ws.Range("J10").Value = 0.4870869789874365
$ws.Range("O10").Value = 0.136569704933115
$ws.Range("P10").Value = 0.136569704933115
$ws.Range("S10").Value = 0.06652132499707661
$ws.Range("T10").Value = 0.0665213249970766
